$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 93.62780766666667
$ws.Range("H2").Value = 280.883423
$ws.Range("I2").Value = 0.3228593149748609
$ws.Range("J2").Value = 0.3228593149748609
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 255.0423041303461
$ws.Range("R2").Value = 2295.380737173115
$ws.Range("S2").Value = 0.01496614507677669
$ws.Range("T2").Value = 0.01496614507677669

$ws.Range("G3").Value = 93.62780766666667
$ws.Range("H3").Value = 280.883423
$ws.Range("I3").Value = 0.3228593149748609
$ws.Range("J3").Value = 0.3228593149748609
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 3803.220595357369
$ws.Range("R3").Value = 34228.98535821632
$ws.Range("S3").Value = 0.2231769015073404
$ws.Range("T3").Value = 0.2231769015073404

$ws.Range("G4").Value = 93.62780766666667
$ws.Range("H4").Value = 280.883423
$ws.Range("I4").Value = 0.3228593149748609
$ws.Range("J4").Value = 0.3228593149748609
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 1443.673850337518
$ws.Range("R4").Value = 12993.06465303767
$ws.Range("S4").Value = 0.08471626839074378
$ws.Range("T4").Value = 0.08471626839074378

$ws.Range("G5").Value = 66.39541
$ws.Range("I5").Value = 0.228953097635189
$ws.Range("J5").Value = 0.228953097635189
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 180.8612074990167
$ws.Range("R5").Value = 1627.75086749115
$ws.Range("S5").Value = 0.01061312192665855
$ws.Range("T5").Value = 0.01061312192665855

$ws.Range("G6").Value = 66.39541
$ws.Range("I6").Value = 0.228953097635189
$ws.Range("J6").Value = 0.228953097635189
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("Q6").Value = 2697.023427571907
$ws.Range("R6").Value = 24273.21084814716
$ws.Range("S6").Value = 0.1582641124190816
$ws.Range("T6").Value = 0.1582641124190816

$ws.Range("G7").Value = 66.39541
$ws.Range("I7").Value = 0.228953097635189
$ws.Range("J7").Value = 0.228953097635189
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 1023.76974948185
$ws.Range("R7").Value = 9213.92774533665
$ws.Range("S7").Value = 0.0600758632894488
$ws.Range("T7").Value = 0.0600758632894488

$ws.Range("G8").Value = 129.9724656666667
$ws.Range("H8").Value = 389.917397
$ws.Range("I8").Value = 0.4481875873899502
$ws.Range("J8").Value = 0.4481875873899502
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 354.0452130967761
$ws.Range("R8").Value = 3186.406917870985
$ws.Range("S8").Value = 0.02077573773893069
$ws.Range("T8").Value = 0.02077573773893069

$ws.Range("G9").Value = 129.9724656666667
$ws.Range("H9").Value = 389.917397
$ws.Range("I9").Value = 0.4481875873899502
$ws.Range("J9").Value = 0.4481875873899502
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 5279.563524681681
$ws.Range("R9").Value = 47516.07172213513
$ws.Range("S9").Value = 0.3098102250992132
$ws.Range("T9").Value = 0.3098102250992132

$ws.Range("G10").Value = 129.9724656666667
$ws.Range("H10").Value = 389.917397
$ws.Range("I10").Value = 0.4481875873899502
$ws.Range("J10").Value = 0.4481875873899502
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 2004.082490267048
$ws.Range("R10").Value = 18036.74241240343
$ws.Range("S10").Value = 0.1176016245518063
$ws.Range("T10").Value = 0.1176016245518063

